# Automatische test-sync: 2025-07-29 21:55:50
# Append a new "Testmail #13" row (row 26) to the historical responses log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

$question = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$reply = "Beste klant,`nBedankt voor je bericht. Helaas kan ik je momenteel niet de datasheet van de VentiQ-250 sturen, aangezien ik niet beschik over het specifieke materiaal waar je naar vraagt.`nIk raad je aan om contact op te nemen met onze verkoopafdeling of de klantenservice, zodat zij je verder kunnen helpen met het verkrijgen van de juiste informatie.`nMet vriendelijke groet,`n[Naam]  `nE-mailassistent bij [Bedrijfsnaam]"
$subject = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$sender = "mailmind.test@zohomail.eu"
$category = "Productinformatie"
$timestamp = "2025-07-29 21:55:30"
$yes = "Ja"
$no = "Nee"

$ws.Cells.Item($row, 1).Value = $question
$ws.Cells.Item($row, 2).Value = $reply
$ws.Cells.Item($row, 3).Value = $subject
$ws.Cells.Item($row, 4).Value = $sender
$ws.Cells.Item($row, 5).Value = $category
$ws.Cells.Item($row, 6).Value = $timestamp
$ws.Cells.Item($row, 7).Value = $yes
$ws.Cells.Item($row, 8).Value = $no
$ws.Cells.Item($row, 9).Value = $yes
$ws.Cells.Item($row, 10).Value = $no
